# Insert a new price record as row 482 in the "Zapallo italiano" sheet.
# This pushes the existing rows 482-551 down to 483-552 (dimension grows
# from A1:R551 to A1:R552), matching the commit's weekly data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 482 (and everything below it) down by one row.
$ws.Rows(482).Insert()

# Populate the newly-opened row 482 with the new record.
$ws.Cells.Item(482, 1).Value = 8
$ws.Cells.Item(482, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(482, 3).Value = "Coquimbo"
$ws.Cells.Item(482, 4).Value = 45127
$ws.Cells.Item(482, 5).Value = 4
$ws.Cells.Item(482, 6).Value = 100112032
$ws.Cells.Item(482, 7).Value = "Zapallo italiano"
$ws.Cells.Item(482, 8).Value = "Sin especificar"
$ws.Cells.Item(482, 9).Value = "Primera"
$ws.Cells.Item(482, 10).Value = 500
$ws.Cells.Item(482, 11).Value = 14000
$ws.Cells.Item(482, 12).Value = 15000
$ws.Cells.Item(482, 13).Value = 14500
$ws.Cells.Item(482, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(482, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(482, 16).Value = 290
$ws.Cells.Item(482, 17).Value = 50
$ws.Cells.Item(482, 18).Value = "Hortaliza"
